$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$values = @(
  "68×87=",
  "38×21=",
  "42×12=",
  "21×48=",
  "55×26=",
  "89×75=",
  "62×57=",
  "93×90=",
  "98×23=",
  "38×71=",
  "76×20=",
  "76×90=",
  "48×57=",
  "53×57=",
  "14×75=",
  "71×65=",
  "15×55=",
  "52×72=",
  "62×19=",
  "57×62=",
  "19×53=",
  "34×83=",
  "37×13=",
  "25×42=",
  "67×49=",
  "66×76=",
  "55×67=",
  "67×99=",
  "17×81=",
  "22×23=",
  "42×65=",
  "18×99=",
  "94×57=",
  "39×92=",
  "73×32=",
  "38×60=",
  "39×10=",
  "39×89=",
  "80×96=",
  "65×50=",
  "95×25=",
  "43×20=",
  "33×100=",
  "29×78=",
  "21×46=",
  "40×97=",
  "89×35=",
  "76×96=",
  "33×22=",
  "66×29=",
  "46×56=",
  "58×62=",
  "88×85=",
  "45×41=",
  "11×94=",
  "37×69=",
  "27×12=",
  "57×46=",
  "14×91=",
  "28×60=",
  "100×73=",
  "89×98=",
  "86×89=",
  "32×28=",
  "75×93=",
  "78×82=",
  "83×56=",
  "79×99=",
  "51×69=",
  "22×98=",
  "81×79=",
  "31×52=",
  "91×23=",
  "35×95=",
  "42×54=",
  "26×87=",
  "38×11=",
  "82×78=",
  "18×45=",
  "21×16=",
  "47×98=",
  "93×92=",
  "47×82=",
  "73×23=",
  "37×42=",
  "49×61=",
  "78×99=",
  "67×59=",
  "87×84=",
  "50×87=",
  "24×33=",
  "27×58=",
  "53×52=",
  "57×36=",
  "83×61=",
  "66×41=",
  "55×79=",
  "34×68=",
  "23×53=",
  "80×42="
)

$numCols = 5
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
  for ($c = 1; $c -le $numCols; $c++) {
    $cell = $t.Cell($r, $c)
    $cell.Range.Text = $values[$idx]
    $idx++
  }
}

Write-Output "done: updated $idx cells"
